# Apply updated loading_percent values (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 13.36918167694865
$ws.Range("C2").Value = 9.272512497948107
$ws.Range("D2").Value = 5.778062008106062
$ws.Range("E2").Value = 12.1838076560925
$ws.Range("F2").Value = 28.04071463641512
$ws.Range("K2").Value = 9.725229834376744
$ws.Range("L2").Value = 9.654118617904535
$ws.Range("N2").Value = 19.7023607037884
$ws.Range("O2").Value = 25.14028807524863
$ws.Range("B3").Value = 13.13119030080535
$ws.Range("C3").Value = 9.272639888474773
$ws.Range("D3").Value = 5.735678529978734
$ws.Range("E3").Value = 12.19356746323869
$ws.Range("F3").Value = 28.06402801069375
$ws.Range("K3").Value = 9.555774066044942
$ws.Range("L3").Value = 9.640379588668608
$ws.Range("N3").Value = 19.76338015506624
$ws.Range("O3").Value = 25.20186646751767
$ws.Range("B4").Value = 12.98559683549393
$ws.Range("C4").Value = 9.272997837982807
$ws.Range("D4").Value = 5.709035498388669
$ws.Range("E4").Value = 12.20178365642184
$ws.Range("F4").Value = 28.08503175340089
$ws.Range("K4").Value = 9.451911282571546
$ws.Range("L4").Value = 9.633608499674478
$ws.Range("N4").Value = 19.80259224988303
$ws.Range("O4").Value = 25.24461826394898
$ws.Range("B5").Value = 12.9264819509243
$ws.Range("C5").Value = 9.27321444502256
$ws.Range("D5").Value = 5.698026013061099
$ws.Range("E5").Value = 12.20569123904664
$ws.Range("F5").Value = 28.09527168203802
$ws.Range("K5").Value = 9.40968817855599
$ws.Range("L5").Value = 9.631270221080007
$ws.Range("N5").Value = 19.81901180687411
$ws.Range("O5").Value = 25.26328090185238
$ws.Range("B6").Value = 12.91668138473353
$ws.Range("C6").Value = 9.273254696618094
$ws.Range("D6").Value = 5.696188805371674
$ws.Range("E6").Value = 12.20637388196929
$ws.Range("F6").Value = 28.09707347280736
$ws.Range("K6").Value = 9.402684851734813
$ws.Range("L6").Value = 9.630907443025398
$ws.Range("N6").Value = 19.82176489333643
$ws.Range("O6").Value = 25.26645471008668
$ws.Range("B7").Value = 12.98479861037407
$ws.Range("C7").Value = 9.273000472270084
$ws.Range("D7").Value = 5.70888763234242
$ws.Range("E7").Value = 12.20183409016055
$ws.Range("F7").Value = 28.08516304976966
$ws.Range("K7").Value = 9.451341360798542
$ws.Range("L7").Value = 9.633575257290675
$ws.Range("N7").Value = 19.80281190526856
$ws.Range("O7").Value = 25.24486493292666
$ws.Range("B8").Value = 13.28705985826047
$ws.Range("C8").Value = 9.272498602780955
$ws.Range("D8").Value = 5.763578839534339
$ws.Range("E8").Value = 12.18671153180527
$ws.Range("F8").Value = 28.04736415674
$ws.Range("K8").Value = 9.666796492801961
$ws.Range("L8").Value = 9.64903739697912
$ws.Range("N8").Value = 19.72303862707965
$ws.Range("O8").Value = 25.1604934234587
$ws.Range("B9").Value = 13.88051993990096
$ws.Range("C9").Value = 9.273716838910687
$ws.Range("D9").Value = 5.865762446356571
$ws.Range("E9").Value = 12.17468212701302
$ws.Range("F9").Value = 28.02634610395997
$ws.Range("K9").Value = 10.08836208218168
$ws.Range("L9").Value = 9.692451669269396
$ws.Range("N9").Value = 19.58039658848083
$ws.Range("O9").Value = 25.0343357762852
$ws.Range("B10").Value = 14.31240521158933
$ws.Range("C10").Value = 9.275932052549802
$ws.Range("D10").Value = 5.937551829611298
$ws.Range("E10").Value = 12.17655675785764
$ws.Range("F10").Value = 28.04327429014515
$ws.Range("K10").Value = 10.39439118618923
$ws.Range("L10").Value = 9.732158788330286
$ws.Range("N10").Value = 19.48392081074432
$ws.Range("O10").Value = 24.96570272352355
$ws.Range("B11").Value = 14.50706629594707
$ws.Range("C11").Value = 9.27722169438376
$ws.Range("D11").Value = 5.969456681710498
$ws.Range("E11").Value = 12.17972496434739
$ws.Range("F11").Value = 28.05798632658008
$ws.Range("K11").Value = 10.53218297331903
$ws.Range("L11").Value = 9.751876467942912
$ws.Range("N11").Value = 19.44182059080064
$ws.Range("O11").Value = 24.93972091132089
$ws.Range("B12").Value = 14.58044749519699
$ws.Range("C12").Value = 9.277750173795015
$ws.Range("D12").Value = 5.981426595320527
$ws.Range("E12").Value = 12.18125636423915
$ws.Range("F12").Value = 28.06456257728848
$ws.Range("K12").Value = 10.58410726971059
$ws.Range("L12").Value = 9.759576965336885
$ws.Range("N12").Value = 19.42613397200694
$ws.Range("O12").Value = 24.93063685954016
$ws.Range("B13").Value = 14.56465946377105
$ws.Range("C13").Value = 9.27763457868064
$ws.Range("D13").Value = 5.978853686726512
$ws.Range("E13").Value = 12.18091182146929
$ws.Range("F13").Value = 28.06310161969438
$ws.Range("K13").Value = 10.57293652942042
$ws.Range("L13").Value = 9.757908193489843
$ws.Range("N13").Value = 19.42950100671261
$ws.Range("O13").Value = 24.93255968993849
$ws.Range("B14").Value = 14.51311048532633
$ws.Range("C14").Value = 9.277264371058235
$ws.Range("D14").Value = 5.970443713684759
$ws.Range("E14").Value = 12.17984431449424
$ws.Range("F14").Value = 28.05850723542914
$ws.Range("K14").Value = 10.53646019345647
$ws.Range("L14").Value = 9.752505330438401
$ws.Range("N14").Value = 19.44052492281626
$ws.Range("O14").Value = 24.93895842936899
$ws.Range("B15").Value = 14.48148980993376
$ws.Range("C15").Value = 9.277042820474753
$ws.Range("D15").Value = 5.965277703422068
$ws.Range("E15").Value = 12.17923358705983
$ws.Range("F15").Value = 28.05582383192399
$ws.Range("K15").Value = 10.51408280609756
$ws.Range("L15").Value = 9.7492262473843
$ws.Range("N15").Value = 19.44731067087346
$ws.Range("O15").Value = 24.94297615874865
$ws.Range("B16").Value = 14.29964098801476
$ws.Range("C16").Value = 9.275853408628123
$ws.Range("D16").Value = 5.935451325052981
$ws.Range("E16").Value = 12.17639619508241
$ws.Range("F16").Value = 28.04245366503965
$ws.Range("K16").Value = 10.38535322340602
$ws.Range("L16").Value = 9.730903143688058
$ws.Range("N16").Value = 19.48670802053663
$ws.Range("O16").Value = 24.96750624263573
$ws.Range("B17").Value = 14.18756701828569
$ws.Range("C17").Value = 9.275195653114395
$ws.Range("D17").Value = 5.916958526178123
$ws.Range("E17").Value = 12.1752477325567
$ws.Range("F17").Value = 28.03604534447198
$ws.Range("K17").Value = 10.30598119744872
$ws.Range("L17").Value = 9.720083487309262
$ws.Range("N17").Value = 19.51133390676499
$ws.Range("O17").Value = 24.98389761409558
$ws.Range("B18").Value = 14.12293922564209
$ws.Range("C18").Value = 9.274843891465267
$ws.Range("D18").Value = 5.906251337122409
$ws.Range("E18").Value = 12.17480519617139
$ws.Range("F18").Value = 28.03301962783171
$ws.Range("K18").Value = 10.26019746299532
$ws.Range("L18").Value = 9.71401631912601
$ws.Range("N18").Value = 19.52566634377319
$ws.Range("O18").Value = 24.99381862299775
$ws.Range("B19").Value = 14.10103117293333
$ws.Range("C19").Value = 9.274729365844131
$ws.Range("D19").Value = 5.902614047137746
$ws.Range("E19").Value = 12.17469284349923
$ws.Range("F19").Value = 28.03210863862351
$ws.Range("K19").Value = 10.24467493115005
$ws.Range("L19").Value = 9.711988996210037
$ws.Range("N19").Value = 19.53054799826046
$ws.Range("O19").Value = 24.99726235991109
$ws.Range("B20").Value = 14.19951518514541
$ws.Range("C20").Value = 9.275262926296769
$ws.Range("D20").Value = 5.918934451062703
$ws.Range("E20").Value = 12.17534743186173
$ws.Range("F20").Value = 28.0366592146169
$ws.Range("K20").Value = 10.31444441711294
$ws.Range("L20").Value = 9.721219140453599
$ws.Range("N20").Value = 19.50869503137073
$ws.Range("O20").Value = 24.98210167835192
$ws.Range("B21").Value = 14.52826126427391
$ws.Range("C21").Value = 9.277372024533227
$ws.Range("D21").Value = 5.97291698778196
$ws.Range("E21").Value = 12.18014887703715
$ws.Range("F21").Value = 28.05982946780632
$ws.Range("K21").Value = 10.54718146924054
$ws.Range("L21").Value = 9.754085969944883
$ws.Range("N21").Value = 19.43727999773177
$ws.Range("O21").Value = 24.93705847347633
$ws.Range("B22").Value = 14.7411451587297
$ws.Range("C22").Value = 9.278984095371746
$ws.Range("D22").Value = 6.007544665485986
$ws.Range("E22").Value = 12.18521937814688
$ws.Range("F22").Value = 28.08082912571273
$ws.Range("K22").Value = 10.69778399166001
$ws.Range("L22").Value = 9.776927255287404
$ws.Range("N22").Value = 19.39209673296483
$ws.Range("O22").Value = 24.91201950397994
$ws.Range("B23").Value = 14.62772852228625
$ws.Range("C23").Value = 9.278102463288658
$ws.Range("D23").Value = 5.989124130666118
$ws.Range("E23").Value = 12.1823368120472
$ws.Range("F23").Value = 28.06908658935565
$ws.Range("K23").Value = 10.61755811234484
$ws.Range("L23").Value = 9.764613331122456
$ws.Range("N23").Value = 19.41607588376609
$ws.Range("O23").Value = 24.92498037905655
$ws.Range("B24").Value = 14.19411402060214
$ws.Range("C24").Value = 9.27523242985718
$ws.Range("D24").Value = 5.918041370104796
$ws.Range("E24").Value = 12.17530167951216
$ws.Range("F24").Value = 28.03637963247326
$ws.Range("K24").Value = 10.31061866418095
$ws.Range("L24").Value = 9.720705234173655
$ws.Range("N24").Value = 19.5098875226939
$ws.Range("O24").Value = 24.98291207164154
$ws.Range("B25").Value = 13.7203768405384
$ws.Range("C25").Value = 9.273153803777564
$ws.Range("D25").Value = 5.838681877528059
$ws.Range("E25").Value = 12.17605202975852
$ws.Range("F25").Value = 28.02634389303926
$ws.Range("K25").Value = 9.97474682825227
$ws.Range("L25").Value = 9.67932182067581
$ws.Range("N25").Value = 19.61751725049814
$ws.Range("O25").Value = 25.06424696071553
